# Update "想去人数" (F column) values on the 展览, 演出, and 全部类型 sheets
# to match the regenerated site output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F7").Value  = 702
$wsExhibition.Range("F9").Value  = 1308
$wsExhibition.Range("F11").Value = 419
$wsExhibition.Range("F14").Value = 50
$wsExhibition.Range("F15").Value = 50
$wsExhibition.Range("F16").Value = 1102
$wsExhibition.Range("F17").Value = 133
$wsExhibition.Range("F19").Value = 421
$wsExhibition.Range("F23").Value = 156
$wsExhibition.Range("F27").Value = 17

# --- Sheet "演出" (Performance) ---
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F7").Value = 250

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value  = 702
$wsAll.Range("F11").Value = 1308
$wsAll.Range("F15").Value = 419
$wsAll.Range("F19").Value = 50
$wsAll.Range("F20").Value = 50
$wsAll.Range("F21").Value = 1102
$wsAll.Range("F23").Value = 133
$wsAll.Range("F25").Value = 421
$wsAll.Range("F28").Value = 250
$wsAll.Range("F35").Value = 156
$wsAll.Range("F39").Value = 17
